$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh - GitHub Actions scheduled update

$ws.Range("D2").Value = '25.921.49'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.589.42'
$ws.Range("E3").Value = '  -1.67%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '209.95'
$ws.Range("E5").Value = '  -1.28%  '
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("D7").Value = '0.484'
$ws.Range("E7").Value = '  -3.24%  '
$ws.Range("D8").Value = '0.248'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").Value = '0.0618'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '18.28'
$ws.Range("E10").Value = '  -0.51%  '
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  -0.38%  '
$ws.Range("D12").Value = '1.808.38'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").Value = '1.583.16'
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").Value = '25.908.72'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").Value = '60.33'
$ws.Range("E17").Value = '  -1.79%  '
$ws.Range("D18").Value = '0.0₃0724'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").Value = '193.98'
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("E21").Value = '  -0.70%  '
$ws.Range("D22").Value = '9.42'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '5.95'
$ws.Range("E23").Value = '  -1.18%  '
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").Value = '141.58'
$ws.Range("E25").Value = '  -1.45%  '
$ws.Range("E26").Value = '  -0.27%  '
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("D28").Value = '15.14'
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("D29").Value = '6.48'
$ws.Range("E29").Value = '  -2.42%  '
$ws.Range("E30").Value = '  -5.47%  '
$ws.Range("D31").Value = '0.0474'
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("E33").Value = '  -1.56%  '
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("D35").Value = '2.35'
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("D36").Value = '1.108.03'
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  -1.44%  '
$ws.Range("D39").Value = '0.506'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("D41").Value = '0.784'
$ws.Range("E41").Value = '  -6.37%  '
$ws.Range("D42").Value = '0.816'
$ws.Range("E42").Value = '  +8.97%  '
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("D44").Value = '93.59'
$ws.Range("E44").Value = '  -4.68%  '
$ws.Range("D45").Value = '1.722.01'
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("E46").Value = '  -1.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.50'
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").Value = '53.51'
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("E51").Value = '  -0.19%  '
